# Logged Week 15 and simulated Week 16
# Update the "H" (home) row on both the OFF and DEF sheets with new
# cumulative target-depth totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (H) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 258
$wsOff.Range("C2").Value = 177
$wsOff.Range("D2").Value = 35
$wsOff.Range("E2").Value = 18
$wsOff.Range("F2").Value = 3
$wsOff.Range("G2").Value = 4

# --- DEF sheet: row 2 (H) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 303
$wsDef.Range("C2").Value = 192
$wsDef.Range("D2").Value = 80
$wsDef.Range("E2").Value = 38
